$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update P_Value_T (column G) values for rows 2-21, and
# Condition_Significant (column F) flips for rows 6 and 16.

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0

$ws.Range("F6").Value = $true
$ws.Range("G6").Value = 0

$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = 0

$ws.Range("G11").Value = [double]"5.050798564436851e-256"

$ws.Range("G12").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("G15").Value = 0

$ws.Range("F16").Value = $true
$ws.Range("G16").Value = [double]"1.982436403432171e-284"

$ws.Range("G17").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("G21").Value = 0
